$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (ownTeam, oppTeam) before the current "batsman" column (D) ---
$ws.Range("D:D").Insert()
$ws.Range("D:D").Insert()

# --- Insert a new row before the current Sharjah row (row 3) ---
$ws.Range("3:3").Insert()

# The new numeric-looking values (runs/balls/4s/6s/sr) must stay text, matching
# the source data (t="str"), otherwise Excel auto-coerces them to numbers.
$ws.Range("G3:K3").NumberFormat = "@"

# --- Header row ---
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# --- Row 2 (Abu Dhabi match) : fill in new ownTeam/oppTeam columns ---
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Sunrisers Hyderabad"

# --- Row 3 (new Dubai (DSC) match) ---
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 05 2020"
$ws.Range("C3").Value = "Capitals won by 59 runs"
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Delhi Capitals"
$ws.Range("F3").Value = "Mohammed Siraj "
$ws.Range("G3").Value = "5"
$ws.Range("H3").Value = "4"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "125.00"

# --- Row 4 (previously Sharjah row 3) : fill in new ownTeam/oppTeam columns ---
$ws.Range("D4").Value = "Royal Challengers Bangalore"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
